$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Append a new row 18 that duplicates the original row 2
# (copy formatting/values before row 2 gets overwritten below)
$ws.Range("A2:R2").Copy($ws.Range("A18:R18"))

# Step 2: Update the weekly price/date figures for rows 2-17
# Row 2
$ws.Range("D2").Value = 44978

# Row 3
$ws.Range("D3").Value = 44985
$ws.Range("J3").Value = 1000

# Row 4
$ws.Range("D4").Value = 44951
$ws.Range("J4").Value = 800
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 2250
$ws.Range("P4").Value = 750

# Row 6
$ws.Range("D6").Value = 44999
$ws.Range("J6").Value = 1100

# Row 7
$ws.Range("D7").Value = 44911
$ws.Range("J7").Value = 700
$ws.Range("K7").Value = 1800
$ws.Range("M7").Value = 1900
$ws.Range("P7").Value = 633

# Row 8
$ws.Range("D8").Value = 44953
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 2000
$ws.Range("L8").Value = 2500
$ws.Range("M8").Value = 2250
$ws.Range("P8").Value = 750

# Row 9
$ws.Range("D9").Value = 44881
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 1900
$ws.Range("L9").Value = 2000
$ws.Range("M9").Value = 1950
$ws.Range("P9").Value = 650

# Row 10
$ws.Range("D10").Value = 44685
$ws.Range("J10").Value = 400
$ws.Range("K10").Value = 1500
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = 1750
$ws.Range("P10").Value = 583

# Row 11
$ws.Range("D11").Value = 44883
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 1800
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 1900
$ws.Range("P11").Value = 633

# Row 12
$ws.Range("D12").Value = 44965
$ws.Range("J12").Value = 1120

# Row 13
$ws.Range("D13").Value = 44970
$ws.Range("J13").Value = 800
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = 2250
$ws.Range("P13").Value = 750

# Row 14
$ws.Range("D14").Value = 44848
$ws.Range("J14").Value = 1000

# Row 15
$ws.Range("D15").Value = 44971
$ws.Range("J15").Value = 1000

# Row 16
$ws.Range("D16").Value = 44992
$ws.Range("J16").Value = 1040
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = 2250
$ws.Range("P16").Value = 750

# Row 17
$ws.Range("D17").Value = 44827
$ws.Range("J17").Value = 1200

